$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows that are no longer part of the new layout (IR0..IR9, old MOTOR/ACK rows)
$ws.Range("A5:B13").Clear()

# ---- Values (write order matters: it controls shared-string insertion order) ----
$ws.Range("A1").Value = "Message Type"
$ws.Range("B1").Value = "ID (Decimal)"
$ws.Range("C1").Value = "Payload Length"

$ws.Range("A2").Value = "SENSOR"
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 20

$ws.Range("A3").Value = "MOTOR"
$ws.Range("B3").Value = 50
$ws.Range("C3").Value = 2

$ws.Range("A4").Value = "ACK"
$ws.Range("B4").Value = 99
$ws.Range("C4").Value = 0

$ws.Range("F1").Value = "Message Index #"
$ws.Range("G1").Value = "Message Type ID"
$ws.Range("H1").Value = "Payload"
$ws.Range("I1").Value = "Terminator"

# ---- Formatting ----
# Bold header for the new C1 cell (A1/B1 already bold from original style)
$ws.Range("C1").Font.Bold = $true

# Colored header cells F1/G1/H1/I1 (fill creation order: I1 red, H1 theme, G1 green, F1 theme)
$ws.Range("I1").Interior.Color = 255

$ws.Range("H1").Interior.Color = 255
$ws.Range("H1").Interior.ThemeColor = 5

$ws.Range("G1").Interior.Color = 5287936

$ws.Range("F1").Interior.Color = 255
$ws.Range("F1").Interior.ThemeColor = 8

# ---- Column widths (closest achievable values that round-trip to the target stored widths) ----
$ws.Columns.Item(1).ColumnWidth = 12.833333333333334
$ws.Columns.Item(2).ColumnWidth = 11.0
$ws.Columns.Item(3).ColumnWidth = 13.666666666666666
$ws.Columns.Item(4).ColumnWidth = 10.0
$ws.Columns.Item(6).ColumnWidth = 14.833333333333334
$ws.Columns.Item(7).ColumnWidth = 15.0
$ws.Columns.Item(8).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 10.0

# ---- Selection ----
[void]$ws.Range("A2").Select()

Write-Host "done"
